$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.181.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.09%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.907.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.57%  "

# Row 4
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.15%  "

# Row 6
$ws.Range("E6").Value = "  +0.11%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5223"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.31%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3766"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.06%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07274"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.76%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.23%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9070"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.05%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08480"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +10.86%  "

# Row 13
$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "97.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.85%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.901.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.45%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.298"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.53%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.13%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008671"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.89%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.63%  "

# Row 19
$ws.Range("E19").Value = "  +0.09%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.218.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.15%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.095"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.34%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.140.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.63%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.44%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.452"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.53%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.322"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.67%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "147.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.48%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.758"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.89%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.31%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.48%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.826"
$ws.Range("D30").Style = "Normal"

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.921"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.55%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09322"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.31%  "

# Row 33
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05067"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.40%  "

# Row 34
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7988"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.04%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.250"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.30%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.446"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.69%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.946"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.83%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.610"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.75%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5723"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.68%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.14%  "

# Row 41
$ws.Range("E41").Value = "  -0.02%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.118"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.03%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.633"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.81%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "115.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.76%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1519"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.57%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4872"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.09%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.35%  "

# Row 48
$ws.Range("E48").Value = "  +0.08%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.629"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.73%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.82"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.36%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.10"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.24%  "
